# Auto-generated PowerShell Excel COM-interop script
# Implements the diff: adds column AN to sheet "data" and column AM to sheet "pocetR"
# with a new survey wave (25. 1. 2022) of data, and updates the "aktualizace" date in the
# footer titles of both sheets from 6. 1. 2022 to 1. 2. 2022.

$wb = $excel.ActiveWorkbook

# ----- Sheet 1: "data" -----
$ws1 = $wb.Worksheets.Item("data")

# Header cell AN1: new wave date, copy style from existing header cell AM1
$ws1.Range("AM1").Copy()
$ws1.Range("AN1").PasteSpecial(-4122)
$ws1.Range("AN1").Value = "25. 1. 2022"

# Data rows 2-76: new percentage values in column AN
$an_values = @{
    2 = 0.58
    3 = 0.26
    4 = 0.16
    5 = 0.71
    6 = 0.16
    7 = 0.13
    8 = 0.55
    9 = 0.28
    10 = 0.17
    11 = 0.76
    12 = 0.14
    13 = 0.1
    14 = 0.66
    15 = 0.22
    16 = 0.12
    17 = 0.59
    18 = 0.2
    19 = 0.21
    20 = 0.4
    21 = 0.36
    22 = 0.24
    23 = 0.76
    24 = 0.13
    25 = 0.11
    26 = 0.64
    27 = 0.24
    28 = 0.12
    29 = 0.41
    30 = 0.34
    31 = 0.25
    32 = 0.39
    33 = 0.3
    34 = 0.31
    35 = 0.51
    36 = 0.31
    37 = 0.18
    38 = 0.67
    39 = 0.23
    40 = 0.1
    41 = 0.72
    42 = 0.18
    43 = 0.1
    44 = 0.61
    45 = 0.23
    46 = 0.16
    47 = 0.43
    48 = 0.35
    49 = 0.22
    50 = 0.63
    51 = 0.24
    52 = 0.13
    53 = 0.52
    54 = 0.31
    55 = 0.17
    56 = 0.59
    57 = 0.23
    58 = 0.18
    59 = 0.62
    60 = 0.23
    61 = 0.15
    62 = 0.55
    63 = 0.28
    64 = 0.17
    65 = 0.62
    66 = 0.22
    67 = 0.16
    68 = 0.54
    69 = 0.29
    70 = 0.17
    71 = 0.54
    72 = 0.3
    73 = 0.16
    74 = 0.55
    75 = 0.28
    76 = 0.17
}
foreach ($row in $an_values.Keys) {
    $ws1.Cells.Item($row, 40).Value = $an_values[$row]
}

# Row 77: update footer title date
$ws1.Range("A77").Value = "Život během pandemie, Počet protektivních aktivit, % respondentů celkově a ve skupinách, aktualizace 1. 2. 2022"

# ----- Sheet 2: "pocetR" -----
$ws2 = $wb.Worksheets.Item("pocetR")

# Header cell AM1: new wave date, copy style from existing header cell AL1
$ws2.Range("AL1").Copy()
$ws2.Range("AM1").PasteSpecial(-4122)
$ws2.Range("AM1").Value = "25. 1. 2022"

# Data rows 2-26: new count values in column AM
$am_values = @{
    2 = 1815
    3 = 412
    4 = 1403
    5 = 304
    6 = 797
    7 = 109
    8 = 605
    9 = 444
    10 = 676
    11 = 695
    12 = 334
    13 = 524
    14 = 725
    15 = 494
    16 = 759
    17 = 562
    18 = 666
    19 = 500
    20 = 649
    21 = 880
    22 = 935
    23 = 952
    24 = 413
    25 = 211
    26 = 239
}
foreach ($row in $am_values.Keys) {
    $ws2.Cells.Item($row, 39).Value = $am_values[$row]
}

# Row 27: update footer title date, and extend the trailing blank cell to AM27
# (copy the existing blank cell AL27 into AM27 so the new cell is created with the
#  same "empty text" content/format as the rest of the row)
$ws2.Range("A27").Value = "Život během pandemie, Počet protektivních aktivit, velikost dotázaného souboru celkově a ve skupinách, aktualizace 1. 2. 2022"
$ws2.Range("AL27").Copy()
$ws2.Range("AM27").PasteSpecial(-4122)

